# New crime data collected - weekly refresh for 122nd Precinct CompStat report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bulletin volume/number and reporting week dates ---
# "Volume 30   Number  48" -> "Volume 30   Number  49"
$ws.Range("A8").Value = "Volume 30   Number  49"
# "Report Covering the Week  11/27/2023  Through  12/3/2023" -> "...12/4/2023  Through  12/10/2023"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Row 16: Robbery ---
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -75
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = -24.193548387096
$ws.Range("L16").Value = 11.904761904761
$ws.Range("M16").Value = -57.272727272727
$ws.Range("N16").Value = -84.067796610169

# --- Row 17: Fel. Assault ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 161
$ws.Range("J17").Value = 122
$ws.Range("K17").Value = 31.967213114754
$ws.Range("L17").Value = 36.440677966101
$ws.Range("M17").Value = 23.846153846153
$ws.Range("N17").Value = -44.482758620689

# --- Row 18: Burglary ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 88.888888888888
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = 22.077922077922
$ws.Range("L18").Value = 40.298507462686
$ws.Range("M18").Value = -52.763819095477
$ws.Range("N18").Value = -92.684824902723

# --- Row 19: Gr. Larceny ---
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -21.875
$ws.Range("I19").Value = 445
$ws.Range("J19").Value = 348
$ws.Range("K19").Value = 27.873563218390
$ws.Range("L19").Value = 39.937106918239
$ws.Range("M19").Value = 14.690721649484
$ws.Range("N19").Value = -44.165621079046

# --- Row 20: G.L.A. ---
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 105
$ws.Range("J20").Value = 130
$ws.Range("K20").Value = -19.230769230769
$ws.Range("L20").Value = 38.157894736842
$ws.Range("M20").Value = -6.25
$ws.Range("N20").Value = -96.043707611153

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = -4.545454545454
$ws.Range("I21").Value = 862
$ws.Range("J21").Value = 748
$ws.Range("K21").Value = 15.240641711229
$ws.Range("L21").Value = 36.825396825396
$ws.Range("M21").Value = -10.114702815432
$ws.Range("N21").Value = -83.863721452639

# --- Row 23: Housing ---
$ws.Range("I16").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1
$ws.Range("I16").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("I21").Copy($ws.Range("E23"))
$ws.Range("E23").Value = 0
$ws.Range("I23").Value = 33
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 65
$ws.Range("L23").Value = -5.714285714285
$ws.Range("M23").Value = 50

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 47.368421052631
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 36.111111111111
$ws.Range("I24").Value = 1043
$ws.Range("J24").Value = 819
$ws.Range("K24").Value = 27.350427350427
$ws.Range("L24").Value = 118.200836820084
$ws.Range("M24").Value = -35.337879727216

# --- Row 25: Misd. Assault ---
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = 23.076923076923
$ws.Range("I25").Value = 311
$ws.Range("J25").Value = 318
$ws.Range("K25").Value = -2.201257861635
$ws.Range("L25").Value = 29.583333333333
$ws.Range("M25").Value = -39.2578125

# --- Row 26: UCR Rape* ---
$ws.Range("I16").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("I21").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = -26.315789473684

# --- Row 27: Other Sex Crimes ---
$ws.Range("I16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = 19.354838709677
$ws.Range("L27").Value = 23.333333333333

# --- Row 30: Hate Crimes ---
$ws.Range("C14").Copy($ws.Range("F30"))
